$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Change 1: paragraph 1 - "elements->printing" becomes "elements and printing",
#     split across 3 runs (Ope | rations: ... and  | printing...) ---
$p1 = $d.Paragraphs.Item(1).Range
$xml1 = '<w:p ' + $wns + '>' +
        '<w:r><w:t>Ope</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">rations: adding 10000 elements and </w:t></w:r>' +
        '<w:r><w:t>printing 10000-&gt;deleting 5000elements-&gt;setting 3995 values to elements</w:t></w:r>' +
        '</w:p>'
$null = $d.Range($p1.Start, $p1.End).InsertXML($xml1)

# --- Change 2: paragraph 3 - "For linked list:  4.44117 seconds" becomes
#     "For linked list: 2.07746 seconds", split across 3 runs ---
$p3 = $d.Paragraphs.Item(3).Range
$xml3 = '<w:p ' + $wns + '>' +
        '<w:r><w:t xml:space="preserve">For linked list: </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">2.07746 </w:t></w:r>' +
        '<w:r><w:t>seconds</w:t></w:r>' +
        '</w:p>'
$null = $d.Range($p3.Start, $p3.End).InsertXML($xml3)

# --- Change 3: merge the "For array: 1.02328" paragraph with the following
#     bookmark-only paragraph (delete the paragraph mark between them), then
#     append a brand-new empty paragraph at the very end of the document. ---
$pArray = $d.Paragraphs.Item(5)
$paraMark = $d.Range($pArray.Range.End - 1, $pArray.Range.End)
$paraMark.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$target = $d.Range($newLast.Range.Start, $newLast.Range.End)
$null = $target.InsertXML('<w:p ' + $wns + '/>')

Write-Output "Done. Paragraphs: $($d.Paragraphs.Count)"
